# Update the "想去人数" (F column) values on both the "展览" and "全部类型"
# sheets, which hold duplicate data sets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1702
    $ws.Range("F3").Value = 7851
    $ws.Range("F4").Value = 182
}
